$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New shared strings must be created in the same order the author typed them
# so they land at the same indices (9..13) in xl/sharedStrings.xml:
#   9  BBC-Eng   10  60x   11  1.10   12  1.20x   13  QQ-eng
# ---------------------------------------------------------------------------

# 9: BBC-Eng
$ws.Range("A14").Value = "BBC-Eng"

# 10: 60x
$ws.Range("A21").Value = "60x"

# 11: 1.10 -- a plain ".Value" assignment of a numeric-looking string like
# "1.10" gets auto-converted to the number 1.1 by Excel's normal text entry
# parsing, and prefixing it with a leading apostrophe marks the cell with a
# "quotePrefix" style that isn't present in the target file. Routing the
# text through a formula + paste-special-values keeps it a genuine string
# without tagging the cell with any extra style.
$ws.Range("H1").Formula = '="1.10"'
$ws.Range("H1").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4163) | Out-Null
$ws.Range("H1").Value = ""

# 12: 1.20x
$ws.Range("C21").Value = "1.20x"

# 13: QQ-eng
$ws.Range("E14").Value = "QQ-eng"

# ---------------------------------------------------------------------------
# New data rows 14-26
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = 13
$ws.Range("C15").Value = 11
$ws.Range("E15").Value = 12
$ws.Range("F15").Value = 13

$ws.Range("A16").Value = 20
$ws.Range("B16").Value = 6
$ws.Range("C16").Value = 10
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 6

$ws.Range("A17").Value = 12
$ws.Range("B17").Value = 24
$ws.Range("C17").Value = 13
$ws.Range("E17").Value = 46
$ws.Range("F17").Value = 80

$ws.Range("A18").Value = 8
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = 16
$ws.Range("E18").Value = 30
$ws.Range("F18").Value = 3

$ws.Range("A19").Value = 45
$ws.Range("B19").Value = 11
$ws.Range("C19").Value = 14
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 5

$ws.Range("A20").Value = 11
$ws.Range("B20").Value = 15
$ws.Range("C20").Value = 14
$ws.Range("E20").Value = 37
$ws.Range("F20").Value = 24

$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 14

$ws.Range("A22").Value = 11
$ws.Range("B22").Value = 50
$ws.Range("C22").Value = 25
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 10

$ws.Range("A23").Value = 10
$ws.Range("B23").Value = 12
$ws.Range("C23").Value = 10
$ws.Range("E23").Value = 5
$ws.Range("F23").Value = 3

$ws.Range("A24").Value = 10
$ws.Range("B24").Value = 7
$ws.Range("C24").Value = 7
$ws.Range("E24").Value = 11
$ws.Range("F24").Value = 6

$ws.Range("A25").Value = 10
$ws.Range("B25").Value = 10
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = 6
$ws.Range("F25").Value = 3

$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 13

# ---------------------------------------------------------------------------
# Final selection left by the author after entering the new block
# ---------------------------------------------------------------------------
$ws.Range("E15:F26").Select() | Out-Null
